$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "bonita" column before the existing "fallbrook" column (BP),
# which shifts fallbrook/lakeside/ramona/ranchosantafe one column to the right.
$ws.Columns("BP:BP").Insert()
$ws.Range("BP1").Value = "bonita"

# Append a new "springvalley" column at the end (after ranchosantafe, now BT).
$ws.Range("BU1").Value = "springvalley"

# Add the new data row for 23-Mar-2020 (row 19).
$ws.Range("A19").Value = 43913
$ws.Range("B19").Value = 213
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 46
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = 32
$ws.Range("I19").Value = 13
$ws.Range("J19").Value = 15
$ws.Range("K19").Value = 6
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 78
$ws.Range("N19").Value = 134
$ws.Range("O19").Value = 1
$ws.Range("P19").Value = 34
$ws.Range("Q19").Value = 17
$ws.Range("R19").Value = 1
$ws.Range("AJ19").Value = 17
$ws.Range("AK19").Value = 0
$ws.Range("AL19").Value = 1
$ws.Range("AM19").Value = 7
$ws.Range("AN19").Value = 3
$ws.Range("AO19").Value = 2
$ws.Range("AP19").Value = 0
$ws.Range("AQ19").Value = 1
$ws.Range("AR19").Value = 3
$ws.Range("AS19").Value = 0
$ws.Range("AT19").Value = 0
$ws.Range("AU19").Value = 6
$ws.Range("AV19").Value = 11
$ws.Range("AW19").Value = 0
$ws.Range("AX19").Value = 3
$ws.Range("AY19").Value = 1
$ws.Range("AZ19").Value = 0
$ws.Range("BA19").Value = 14
$ws.Range("BB19").Value = 6
$ws.Range("BC19").Value = 5
$ws.Range("BD19").Value = 7
$ws.Range("BE19").Value = 5
$ws.Range("BF19").Value = 5
$ws.Range("BG19").Value = 1
$ws.Range("BH19").Value = 3
$ws.Range("BI19").Value = 5
$ws.Range("BJ19").Value = 2
$ws.Range("BK19").Value = 134
$ws.Range("BL19").Value = 3
$ws.Range("BM19").Value = 1
$ws.Range("BN19").Value = 1
$ws.Range("BO19").Value = 2
$ws.Range("BP19").Value = 1
$ws.Range("BQ19").Value = 2
$ws.Range("BR19").Value = 2
$ws.Range("BS19").Value = 2
$ws.Range("BT19").Value = 5
$ws.Range("BU19").Value = 2

# Keep the (stale) hidden _FilterDatabase name in sync with the shifted last column.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$BT`$17"

# Update the active selection to match where the new data was entered.
$null = $ws.Range("BU19").Select()
